$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Historical Account of the Project" paragraph: split the single summary
#    sentence into a bulleted-style list of historical events, re-using the
#    tail of the original sentence in a final paragraph.
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute(
    "This should outline the main events over the lifetime of the project, and how the project team acted to produce a plan and to deliver a product within a constrained lifetime. This should take no more than two pages of A4.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = ""
$xml = @"
<w:p><w:r><w:t>This should outline the main events over the lifetime of the project</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>First meeting get to know each other, interests, what we are good at what we are bad at, etc.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Different people assigned different things to research to help with deciding the different options, such as glassfish and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>google</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> app engine, etc.</w:t></w:r><w:r><w:t xml:space="preserve"> and a project plan is written</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>Initial roles assigned and split into teams, testing and design.</w:t></w:r></w:p><w:p><w:r><w:t>Testing team assigned to creating test specification, design team assigned to creating design specification.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Team members with most experience in programming assigned to carrying out </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>spikework</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> on different things</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="6" w:name="_GoBack"/><w:bookmarkEnd w:id="6"/></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>and</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> how the project team acted to produce a plan and to deliver a product within a constrained lifetime. This should take no more than two pages of A4.</w:t></w:r></w:p>
"@
$null = $r.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2. Collapse the run-split paragraphs below back into single runs (their
#    text doesn't change, only the run boundaries do).
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "This should give a summary of which parts of the project are perceived as correct and which are not. It is as well to be as accurate as possible here - more marks will be deducted for problems that are not declared but are detected by the markers than for problems that are declared in the final report. As well as missing or erroneous features in the software, known problems with documents should be included here.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This should give a summary of which parts of the project are perceived as correct and which are not. It is as well to be as accurate as possible here - more marks will be deducted for problems that are not declared but are detected by the markers than for problems that are declared in the final report. As well as missing or erroneous features in the software, known problems with documents should be included here.",
    2)

$null = $d.Content.Find.Execute(
    "The project leader should write a half page description of the duties and performance of each group member, including the group leaders themselves. This should be agreed with the group member if possible, and it should state whether agreement was reached, and if not, should give an explanation why not.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The project leader should write a half page description of the duties and performance of each group member, including the group leaders themselves. This should be agreed with the group member if possible, and it should state whether agreement was reached, and if not, should give an explanation why not.",
    2)

$null = $d.Content.Find.Execute(
    "This should be no more than a page in length and should address the following subjects:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This should be no more than a page in length and should address the following subjects:",
    2)

# ---------------------------------------------------------------------------
# 3. Remove the stray _GoBack bookmark that used to sit after item 3 in the
#    "Critical Evaluation" list (it now lives in the historical account
#    section instead - see step 1).
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $ptext = $p.Range.Text
    if ($ptext -like "3. What were the most important lessons*") {
        $bms = $p.Range.Bookmarks
        for ($j = $bms.Count; $j -ge 1; $j--) {
            $bm = $bms.Item($j)
            if ($bm.Name -eq "_GoBack") {
                $bm.Delete()
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 4. Footer page-number field cache: "Page 4 of 4" -> "Page 5 of 4" (only the
#    first cached field result changes, matching the authoritative edit).
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$ffields = $footer.Range.Fields
for ($k = 1; $k -le $ffields.Count; $k++) {
    $fld = $ffields.Item($k)
    if ($fld.Type -eq 33) {
        $fld.Result.Text = "5"
        break
    }
}
